$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.536.56'
$ws.Range("E2").Value = '  +4.50%  '
$ws.Range("D3").Value = '3.109.99'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '619.33'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.383'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.80%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.903'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.93%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '3.103.34'
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.676'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +16.10%  '
$ws.Range("E12").Value = '  +6.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000259'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.87%  '
$ws.Range("D14").Value = '91.147.81'
$ws.Range("E14").Value = '  +4.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.39'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '33.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.63%  '
$ws.Range("D17").Value = '3.672.24'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").Value = '3.099.66'
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000227'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '434.61'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '84.17'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.12%  '
$ws.Range("D28").Value = '3.250.44'
$ws.Range("E28").Value = '  -0.53%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.167'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.29%  '
$ws.Range("E31").Value = '  +1.45%  '
$ws.Range("E32").Value = '  +7.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.92'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '519.42'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.03'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.141'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.30'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.11%  '
$ws.Range("E38").Value = '  +2.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.05'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.20%  '
$ws.Range("E40").Value = '  +0.67%  '
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.143'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.372'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.20%  '
$ws.Range("E45").Value = '  +3.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0723'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +10.57%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.89'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '142.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000266'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +17.23%  '
$ws.Range("E50").Value = '  +7.12%  '
$ws.Range("E51").Value = '  +6.07%  '
